$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.713.40"
$ws.Range("E2").Value = "  -0.24%  "
$ws.Range("D3").Value = "1.637.41"
$ws.Range("E3").Value = "  -0.72%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "217.57"
$ws.Range("E5").Value = "  +0.53%  "
$ws.Range("E6").Value = "  -0.93%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("E8").Value = "  -0.54%  "
$ws.Range("D9").Value = "0.0623"
$ws.Range("E9").Value = "  -0.76%  "
$ws.Range("E10").Value = "  -0.63%  "
$ws.Range("D11").Value = "0.0844"
$ws.Range("E11").Value = "  +0.19%  "
$ws.Range("D12").Value = "1.863.77"
$ws.Range("D13").Value = "1.637.21"
$ws.Range("E13").Value = "  -0.82%  "
$ws.Range("E14").Value = "  -1.23%  "
$ws.Range("E15").Value = "  -1.54%  "
$ws.Range("D16").Value = "64.42"
$ws.Range("E16").Value = "  -1.61%  "
$ws.Range("D17").Value = "26.694.18"
$ws.Range("E17").Value = "  -0.34%  "
$ws.Range("E18").Value = "  -2.48%  "
$ws.Range("D19").Value = "211.36"
$ws.Range("E19").Value = "  -3.35%  "
$ws.Range("E20").Value = "  -0.12%  "
$ws.Range("E21").Value = "  -0.92%  "
$ws.Range("E22").Value = "  -1.42%  "
$ws.Range("D23").Value = "2.29"
$ws.Range("E23").Value = "  -2.84%  "
$ws.Range("E24").Value = "  -2.65%  "
$ws.Range("D25").Value = "146.64"
$ws.Range("E25").Value = "  -0.09%  "
$ws.Range("E26").Value = "  -0.25%  "
$ws.Range("D27").Value = "0.118"
$ws.Range("E27").Value = "  -1.98%  "
$ws.Range("E28").Value = "  -0.58%  "
$ws.Range("E29").Value = "  -1.31%  "
$ws.Range("E30").Value = "  -2.59%  "
$ws.Range("E31").Value = "  +0.41%  "
$ws.Range("E32").Value = "  -0.21%  "
$ws.Range("E33").Value = "  -1.36%  "
$ws.Range("D34").Value = "1.272.66"
$ws.Range("E34").Value = "  -0.71%  "
$ws.Range("D35").Value = "1.53"
$ws.Range("E35").Value = "  -1.28%  "
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("E37").Value = "  -1.99%  "
$ws.Range("E38").Value = "  -1.99%  "
$ws.Range("D39").Value = "0.805"
$ws.Range("E39").Value = "  -2.88%  "
$ws.Range("E40").Value = "  -0.11%  "
$ws.Range("E41").Value = "  -1.58%  "
$ws.Range("E42").Value = "  -2.65%  "
$ws.Range("D43").Value = "5.27"
$ws.Range("E43").Value = "  -3.66%  "
$ws.Range("D44").Value = "1.774.13"
$ws.Range("E44").Value = "  -0.82%  "
$ws.Range("D45").Value = "91.47"
$ws.Range("E45").Value = "  -0.64%  "
$ws.Range("D46").Value = "60.27"
$ws.Range("E46").Value = "  +0.80%  "
$ws.Range("E47").Value = "  -1.42%  "
$ws.Range("E48").Value = "  +0.35%  "
$ws.Range("D49").Value = "7.56"
$ws.Range("E49").Value = "  -3.01%  "
$ws.Range("E50").Value = "  -0.95%  "
$ws.Range("E51").Value = "  -0.54%  "
